$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: split the run containing
#   "возобновляемые и невозобновляемые. Только совместными усилиями многих "
# into three runs, isolating "невозобновляемые" (which a spell-checker would
# flag) from its neighbours, matching the diff's three-run layout.
# ---------------------------------------------------------------------------
$rng1 = $d.Content.Duplicate
$found1 = $rng1.Find.Execute("возобновляемые и невозобновляемые. Только совместными усилиями многих ")
if ($found1) {
    $s = $rng1.Start
    $e = $rng1.End

    $p1Len = 17   # "возобновляемые и "
    $p2Len = 16   # "невозобновляемые"

    # Force a run boundary right after "возобновляемые и " by toggling a
    # character property on and back off (this leaves that run's content
    # untouched but prevents it from being coalesced with its neighbour).
    $left = $d.Range($s, $s + $p1Len)
    $left.Bold = 1
    $left.Bold = 0

    # Force a run boundary right before ". Только совместными усилиями..." so
    # that "невозобновляемые" ends up alone, with no direct formatting touch.
    $right = $d.Range($s + $p1Len + $p2Len, $e)
    $right.Bold = 1
    $right.Bold = 0
}

# ---------------------------------------------------------------------------
# Hunk 2: "Also, environmental protection should be raised ..." becomes
#   "Also, " + "ecological" + " " + "protection should be raised ..."
# (word swap "environmental" -> "ecological", plus a four-way run split).
# ---------------------------------------------------------------------------
$rng2 = $d.Content.Duplicate
$found2 = $rng2.Find.Execute("Also, environmental protection should be raised at the global community level. International programs for the study of natural resources should be developed.")
if ($found2) {
    $s2 = $rng2.Start

    $lead = "Also, "
    $oldWord = "environmental"
    $newWord = "ecological"

    $leadLen = $lead.Length
    $oldWordLen = $oldWord.Length

    # Replace "environmental" with "ecological" first.
    $wordRange = $d.Range($s2 + $leadLen, $s2 + $leadLen + $oldWordLen)
    $wordRange.Text = $newWord
    $newWordLen = $newWord.Length

    # Split "Also, " away from "ecological".
    $r1 = $d.Range($s2, $s2 + $leadLen)
    $r1.Bold = 1
    $r1.Bold = 0

    # Split "ecological" away from the following space.
    $r2 = $d.Range($s2 + $leadLen, $s2 + $leadLen + $newWordLen)
    $r2.Bold = 1
    $r2.Bold = 0

    # Split the single space away from "protection should be raised...".
    $r3 = $d.Range($s2 + $leadLen + $newWordLen, $s2 + $leadLen + $newWordLen + 1)
    $r3.Bold = 1
    $r3.Bold = 0
}

Write-Host "done"
